$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Total" column header ---
$ws.Range("X1").Value = "Total"

# --- Row totals for existing category rows (2-6), column X = sum(B:W) ---
$ws.Range("X2").Value = 2269
$ws.Range("X3").Value = 278
$ws.Range("X4").Value = 878
$ws.Range("X5").Value = 301
$ws.Range("X6").Value = 1633

# --- New category row 7: "Outros" ---
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 125
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 39
$ws.Range("F7").Value = 76
$ws.Range("G7").Value = 78
$ws.Range("H7").Value = 68
$ws.Range("I7").Value = 84
$ws.Range("J7").Value = 78
$ws.Range("K7").Value = 91
$ws.Range("L7").Value = 114
$ws.Range("M7").Value = 99
$ws.Range("N7").Value = 126
$ws.Range("O7").Value = 135
$ws.Range("P7").Value = 137
$ws.Range("Q7").Value = 159
$ws.Range("R7").Value = 232
$ws.Range("S7").Value = 225
$ws.Range("T7").Value = 167
$ws.Range("U7").Value = 59
$ws.Range("V7").Value = 16
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 2119

# --- New grand-total row 8: "Total" ---
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 139
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 13
$ws.Range("E8").Value = 42
$ws.Range("F8").Value = 89
$ws.Range("G8").Value = 90
$ws.Range("H8").Value = 104
$ws.Range("I8").Value = 138
$ws.Range("J8").Value = 192
$ws.Range("K8").Value = 244
$ws.Range("L8").Value = 347
$ws.Range("M8").Value = 486
$ws.Range("N8").Value = 600
$ws.Range("O8").Value = 662
$ws.Range("P8").Value = 738
$ws.Range("Q8").Value = 870
$ws.Range("R8").Value = 956
$ws.Range("S8").Value = 868
$ws.Range("T8").Value = 604
$ws.Range("U8").Value = 236
$ws.Range("V8").Value = 50
$ws.Range("W8").Value = 1
$ws.Range("X8").Value = 7478
